$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 194
$ws.Range("F5").Value = 176
$ws.Range("F6").Value = 832
$ws.Range("F7").Value = 4237
$ws.Range("F11").Value = 6175
$ws.Range("F12").Value = 6175
$ws.Range("F13").Value = 69
$ws.Range("F15").Value = 2362
$ws.Range("F17").Value = 170
$ws.Range("F18").Value = 486
$ws.Range("F19").Value = 9306
$ws.Range("F20").Value = 46
$ws.Range("F21").Value = 2511
$ws.Range("F23").Value = 2329
$ws.Range("F24").Value = 2477
$ws.Range("F26").Value = 248
$ws.Range("F27").Value = 1981
$ws.Range("F30").Value = 337
$ws.Range("F32").Value = 48
$ws.Range("F34").Value = 45
$ws.Range("F35").Value = 80
$ws.Range("F36").Value = 586
$ws.Range("F37").Value = 1227
$ws.Range("F42").Value = 1566
$ws.Range("F43").Value = 2578
$ws.Range("F46").Value = 313
$ws.Range("F47").Value = 1256
$ws.Range("F50").Value = 2

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 152
$ws.Range("F22").Value = 107

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 699
$ws.Range("F3").Value = 910

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 194
$ws.Range("F3").Value = 699
$ws.Range("F4").Value = 910
$ws.Range("F10").Value = 176
$ws.Range("F11").Value = 832
$ws.Range("F12").Value = 4237
$ws.Range("F17").Value = 6175
$ws.Range("F18").Value = 69
$ws.Range("F19").Value = 2362
$ws.Range("F20").Value = 170
$ws.Range("F21").Value = 486
$ws.Range("F22").Value = 9306
$ws.Range("F23").Value = 152
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 2511
$ws.Range("F27").Value = 2477
$ws.Range("F29").Value = 248
$ws.Range("F30").Value = 1981
$ws.Range("F33").Value = 337
$ws.Range("F34").Value = 48
$ws.Range("F36").Value = 45
$ws.Range("F37").Value = 80
$ws.Range("F38").Value = 586
$ws.Range("F39").Value = 1227
$ws.Range("F43").Value = 2578
$ws.Range("F45").Value = 313
$ws.Range("F50").Value = 107
$ws.Range("F51").Value = 107
